$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A210").Value = 45966
$ws.Range("B210").Value = 1827
$ws.Range("C210").Value = 134865
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 29135091
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 53229811257
$ws.Range("H210").Value = 0
$ws.Range("A211").Value = 45967
$ws.Range("B211").Value = 1827
$ws.Range("C211").Value = 36155
$ws.Range("D211").Value = 0
$ws.Range("E211").Value = 29135091
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 53229811257
$ws.Range("H211").Value = 0
$ws.Range("A212").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A212").Value = 45968
$ws.Range("B212").Value = 1812
$ws.Range("C212").Value = 65409
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 29135091
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 52792784892
$ws.Range("H212").Value = 0
$ws.Range("A213").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A213").Value = 45971
$ws.Range("B213").Value = 1823
$ws.Range("C213").Value = 60923
$ws.Range("D213").Value = 0
$ws.Range("E213").Value = 29135091
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 53113270893
$ws.Range("H213").Value = 0
$ws.Range("A214").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A214").Value = 45972
$ws.Range("B214").Value = 1823
$ws.Range("C214").Value = 38963
$ws.Range("D214").Value = 0
$ws.Range("E214").Value = 29135091
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 53113270893
$ws.Range("H214").Value = 0
$ws.Range("A215").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A215").Value = 45973
$ws.Range("B215").Value = 1823
$ws.Range("C215").Value = 59251
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 29135091
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 53113270893
$ws.Range("H215").Value = 0
$ws.Range("A216").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A216").Value = 45974
$ws.Range("B216").Value = 1826
$ws.Range("C216").Value = 82802
$ws.Range("A217").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A217").Value = 45975
$ws.Range("B217").Value = 1818
$ws.Range("C217").Value = 38935

$ws = $wb.Worksheets.Item(2)
$ws.Range("A210").Value = 45966
$ws.Range("B210").Value = 3710
$ws.Range("C210").Value = 1035084
$ws.Range("D210").Value = 212112
$ws.Range("E210").Value = 47816422
$ws.Range("F210").Value = 786935520
$ws.Range("G210").Value = 177398925620
$ws.Range("H210").Value = 0.4399999976158142
$ws.Range("A211").Value = 45967
$ws.Range("B211").Value = 3705
$ws.Range("C211").Value = 1320232
$ws.Range("D211").Value = 208937
$ws.Range("E211").Value = 47816422
$ws.Range("F211").Value = 774111585
$ws.Range("G211").Value = 177159843510
$ws.Range("H211").Value = 0.4399999976158142
$ws.Range("A212").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A212").Value = 45968
$ws.Range("B212").Value = 3625
$ws.Range("C212").Value = 579914
$ws.Range("D212").Value = 237662
$ws.Range("E212").Value = 47816422
$ws.Range("F212").Value = 861524750
$ws.Range("G212").Value = 173334529750
$ws.Range("H212").Value = 0.5
$ws.Range("A213").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A213").Value = 45971
$ws.Range("B213").Value = 3660
$ws.Range("C213").Value = 347868
$ws.Range("D213").Value = 240353
$ws.Range("E213").Value = 47816422
$ws.Range("F213").Value = 879691980
$ws.Range("G213").Value = 175008104520
$ws.Range("H213").Value = 0.5
$ws.Range("A214").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A214").Value = 45972
$ws.Range("B214").Value = 3630
$ws.Range("C214").Value = 364150
$ws.Range("D214").Value = 245890
$ws.Range("E214").Value = 47816422
$ws.Range("F214").Value = 892580700
$ws.Range("G214").Value = 173573611860
$ws.Range("H214").Value = 0.5099999904632568
$ws.Range("A215").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A215").Value = 45973
$ws.Range("B215").Value = 3640
$ws.Range("C215").Value = 325646
$ws.Range("D215").Value = 244330
$ws.Range("E215").Value = 47816422
$ws.Range("F215").Value = 889361200
$ws.Range("G215").Value = 174051776080
$ws.Range("H215").Value = 0.5099999904632568
$ws.Range("A216").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A216").Value = 45974
$ws.Range("B216").Value = 3615
$ws.Range("C216").Value = 316320
$ws.Range("A217").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A217").Value = 45975
$ws.Range("B217").Value = 3510
$ws.Range("C217").Value = 398939

$ws = $wb.Worksheets.Item(3)
$ws.Range("A210").Value = 45966
$ws.Range("B210").Value = 7810
$ws.Range("C210").Value = 59635
$ws.Range("D210").Value = 17825
$ws.Range("E210").Value = 12578946
$ws.Range("F210").Value = 140817500
$ws.Range("G210").Value = 99373673400
$ws.Range("H210").Value = 0.1400000005960464
$ws.Range("A211").Value = 45967
$ws.Range("B211").Value = 7810
$ws.Range("C211").Value = 46205
$ws.Range("D211").Value = 17825
$ws.Range("E211").Value = 12578946
$ws.Range("F211").Value = 139213250
$ws.Range("G211").Value = 98241568260
$ws.Range("H211").Value = 0.1400000005960464
$ws.Range("A212").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A212").Value = 45968
$ws.Range("B212").Value = 7740
$ws.Range("C212").Value = 61682
$ws.Range("D212").Value = 39820
$ws.Range("E212").Value = 12578946
$ws.Range("F212").Value = 308206800
$ws.Range("G212").Value = 97361042040
$ws.Range("H212").Value = 0.3199999928474426
$ws.Range("A213").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A213").Value = 45971
$ws.Range("B213").Value = 8120
$ws.Range("C213").Value = 40906
$ws.Range("D213").Value = 39820
$ws.Range("E213").Value = 12578946
$ws.Range("F213").Value = 323338400
$ws.Range("G213").Value = 102141041520
$ws.Range("H213").Value = 0.3199999928474426
$ws.Range("A214").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A214").Value = 45972
$ws.Range("B214").Value = 8200
$ws.Range("C214").Value = 46371
$ws.Range("D214").Value = 39820
$ws.Range("E214").Value = 12578946
$ws.Range("F214").Value = 326524000
$ws.Range("G214").Value = 103147357200
$ws.Range("H214").Value = 0.3199999928474426
$ws.Range("A215").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A215").Value = 45973
$ws.Range("B215").Value = 8190
$ws.Range("C215").Value = 24622
$ws.Range("D215").Value = 39820
$ws.Range("E215").Value = 12578946
$ws.Range("F215").Value = 326125800
$ws.Range("G215").Value = 103021567740
$ws.Range("H215").Value = 0.3199999928474426
$ws.Range("A216").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A216").Value = 45974
$ws.Range("B216").Value = 8670
$ws.Range("C216").Value = 174497
$ws.Range("A217").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A217").Value = 45975
$ws.Range("B217").Value = 8540
$ws.Range("C217").Value = 53392

$ws = $wb.Worksheets.Item(4)
$ws.Range("A210").Value = 45966
$ws.Range("B210").Value = 2115
$ws.Range("C210").Value = 1133814
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 70922823
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 140143498248
$ws.Range("H210").Value = 0
$ws.Range("A211").Value = 45967
$ws.Range("B211").Value = 2115
$ws.Range("C211").Value = 1021468
$ws.Range("D211").Value = 0
$ws.Range("E211").Value = 70922823
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 150001770645
$ws.Range("H211").Value = 0
$ws.Range("A212").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A212").Value = 45968
$ws.Range("B212").Value = 2120
$ws.Range("C212").Value = 794528
$ws.Range("D212").Value = 0
$ws.Range("E212").Value = 70922823
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 150356384760
$ws.Range("H212").Value = 0
$ws.Range("A213").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A213").Value = 45971
$ws.Range("B213").Value = 2100
$ws.Range("C213").Value = 993438
$ws.Range("D213").Value = 51640
$ws.Range("E213").Value = 70922823
$ws.Range("F213").Value = 108444000
$ws.Range("G213").Value = 148937928300
$ws.Range("H213").Value = 0.07000000029802322
$ws.Range("A214").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A214").Value = 45972
$ws.Range("B214").Value = 2085
$ws.Range("C214").Value = 1022961
$ws.Range("D214").Value = 67064
$ws.Range("E214").Value = 70922823
$ws.Range("F214").Value = 139828440
$ws.Range("G214").Value = 147874085955
$ws.Range("H214").Value = 0.09000000357627869
$ws.Range("A215").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A215").Value = 45973
$ws.Range("B215").Value = 2130
$ws.Range("C215").Value = 950878
$ws.Range("D215").Value = 75146
$ws.Range("E215").Value = 70922823
$ws.Range("F215").Value = 160060980
$ws.Range("G215").Value = 151065612990
$ws.Range("H215").Value = 0.1099999994039536
$ws.Range("A216").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A216").Value = 45974
$ws.Range("B216").Value = 2210
$ws.Range("C216").Value = 1144056
$ws.Range("A217").NumberFormat = $ws.Range("A211").NumberFormat
$ws.Range("A217").Value = 45975
$ws.Range("B217").Value = 2240
$ws.Range("C217").Value = 822467

